# Apply the "cryptos" price-refresh edit described in the commit diff.
# Column D holds price strings that must stay literal TEXT (e.g. "1.004",
# "0.000008794") even though they look numeric -- Excel auto-converts a
# plain assignment of a parseable numeric string into a real number, which
# would silently drop trailing zeros / reformat small decimals. We force
# those specific cells to Text format first so the literal string sticks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.709.43'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.890.96'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.88%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.15'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4813'
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3789'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07337'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9189'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.46'
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07699'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '1.922.32'
$ws.Range("E13").Value = '  +2.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.469'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.593'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008794'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '27.765.21'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.50'
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.125'
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").Value = '2.136.26'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.81'
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.29'
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.901'
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.41'
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.122'
$ws.Range("E28").Value = '  +5.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.38'
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08944'
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.158'
$ws.Range("E32").Value = '  -5.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.234'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7607'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.630'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02032'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.531'
$ws.Range("E37").Value = '  -5.60%  '
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05255'
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5440'
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.973'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.957'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1516'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.307'
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.65'
$ws.Range("E45").Value = '  +4.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.58'
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4779'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.637'
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.66'
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06060'
$ws.Range("E51").Value = '  -0.63%  '
